# Update the "Configuration" sheet's mapping/path templates to replace the
# old "SubjectId" placeholder with "ID", split the old combined filename
# template into separate YYYY/MM/DD, Label and filename extension templates,
# and move the active selection to C9.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Configuration")

# Row 2 ("Paths"): {SubjectName}-{SubjectId} / {YYYY}{MM}{DD}-{Label}.dcm
# becomes four separate path templates.
$ws.Range("C2").Value = "{SubjectName}-{ID}"
$ws.Range("D2").Value = "{YYYY}{MM}{DD}"
$ws.Range("E2").Value = "{Label}"
$ws.Range("F2").Value = "{filename}.dcm"

# Row 7 ("Mappings" / "Subject"): SubjectId -> ID
$ws.Range("C7").Value = "ID"

# Move the active selection to C9, as in the saved workbook.
$ws.Range("C9").Select()
